{"js": "const replacements = [\n  [\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"],\n  [\"210\u00f79=23, 3\", \"381\u00f78=47, 5\"],\n  [\"799\u00f77=114, 1\", \"349\u00f79=38, 7\"],\n  [\"621\u00f75=124, 1\", \"431\u00f79=47, 8\"],\n  [\"664\u00f73=221, 1\", \"104\u00f77=14, 6\"],\n  [\"795\u00f75=159, 0\", \"191\u00f78=23, 7\"],\n  [\"120\u00f73=40, 0\", \"379\u00f77=54, 1\"],\n  [\"802\u00f77=114, 4\", \"126\u00f75=25, 1\"],\n  [\"502\u00f75=100, 2\", \"731\u00f75=146, 1\"],\n  [\"950\u00f76=158, 2\", \"586\u00f79=65, 1\"],\n  [\"557\u00f75=111, 2\", \"589\u00f75=117, 4\"],\n  [\"135\u00f79=15, 0\", \"471\u00f74=117, 3\"],\n  [\"932\u00f79=103, 5\", \"408\u00f74=102, 0\"],\n  [\"881\u00f74=220, 1\", \"313\u00f73=104, 1\"],\n  [\"591\u00f74=147, 3\", \"616\u00f77=88, 0\"],\n  [\"394\u00f77=56, 2\", \"165\u00f75=33, 0\"],\n  [\"271\u00f75=54, 1\", \"700\u00f74=175, 0\"],\n  [\"853\u00f76=142, 1\", \"633\u00f79=70, 3\"],\n  [\"652\u00f78=81, 4\", \"592\u00f79=65, 7\"],\n  [\"167\u00f77=23, 6\", \"625\u00f79=69, 4\"],\n  [\"548\u00f74=137, 0\", \"500\u00f79=55, 5\"],\n  [\"327\u00f72=163, 1\", \"581\u00f74=145, 1\"],\n  [\"280\u00f74=70, 0\", \"530\u00f72=265, 0\"],\n  [\"949\u00f76=158, 1\", \"316\u00f76=52, 4\"],\n  [\"116\u00f79=12, 8\", \"144\u00f74=36, 0\"],\n  [\"357\u00f75=71, 2\", \"219\u00f74=54, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"),\n    @(\"210\u00f79=23, 3\", \"381\u00f78=47, 5\"),\n    @(\"799\u00f77=114, 1\", \"349\u00f79=38, 7\"),\n    @(\"621\u00f75=124, 1\", \"431\u00f79=47, 8\"),\n    @(\"664\u00f73=221, 1\", \"104\u00f77=14, 6\"),\n    @(\"795\u00f75=159, 0\", \"191\u00f78=23, 7\"),\n    @(\"120\u00f73=40, 0\", \"379\u00f77=54, 1\"),\n    @(\"802\u00f77=114, 4\", \"126\u00f75=25, 1\"),\n    @(\"502\u00f75=100, 2\", \"731\u00f75=146, 1\"),\n    @(\"950\u00f76=158, 2\", \"586\u00f79=65, 1\"),\n    @(\"557\u00f75=111, 2\", \"589\u00f75=117, 4\"),\n    @(\"135\u00f79=15, 0\", \"471\u00f74=117, 3\"),\n    @(\"932\u00f79=103, 5\", \"408\u00f74=102, 0\"),\n    @(\"881\u00f74=220, 1\", \"313\u00f73=104, 1\"),\n    @(\"591\u00f74=147, 3\", \"616\u00f77=88, 0\"),\n    @(\"394\u00f77=56, 2\", \"165\u00f75=33, 0\"),\n    @(\"271\u00f75=54, 1\", \"700\u00f74=175, 0\"),\n    @(\"853\u00f76=142, 1\", \"633\u00f79=70, 3\"),\n    @(\"652\u00f78=81, 4\", \"592\u00f79=65, 7\"),\n    @(\"167\u00f77=23, 6\", \"625\u00f79=69, 4\"),\n    @(\"548\u00f74=137, 0\", \"500\u00f79=55, 5\"),\n    @(\"327\u00f72=163, 1\", \"581\u00f74=145, 1\"),\n    @(\"280\u00f74=70, 0\", \"530\u00f72=265, 0\"),\n    @(\"949\u00f76=158, 1\", \"316\u00f76=52, 4\"),\n    @(\"116\u00f79=12, 8\", \"144\u00f74=36, 0\"),\n    @(\"357\u00f75=71, 2\", \"219\u00f74=54, 3\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n    $find.Parent.Text = $newText\n}"}
